$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - rows 2-7, column F changes
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 14813
$wsExhibition.Range("F3").Value = 343
$wsExhibition.Range("F4").Value = 704
$wsExhibition.Range("F5").Value = 245
$wsExhibition.Range("F6").Value = 618
$wsExhibition.Range("F7").Value = 1578

# Sheet "全部类型" (all types) - rows 2,3,4,5,8,9, column F changes
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 14813
$wsAll.Range("F3").Value = 343
$wsAll.Range("F4").Value = 704
$wsAll.Range("F5").Value = 245
$wsAll.Range("F8").Value = 618
$wsAll.Range("F9").Value = 1578
